$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number format on Price (D) column cells being updated so Excel
# does not auto-convert values such as "6.41" or "1.00" into numbers/dates.
$priceCells = @("D2","D3","D5","D6","D9","D10","D11","D13","D14","D17","D18","D19","D20","D21","D22","D23","D25","D26","D28","D29","D30","D32","D33","D34","D35","D37","D38","D40","D41","D42","D44","D45","D46","D48","D49","D50","D51")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated values
$ws.Range("D2").Value = "60.495.09"
$ws.Range("E2").Value = "  +0.13%  "
$ws.Range("D3").Value = "2.628.37"
$ws.Range("E3").Value = "  +0.98%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "523.03"
$ws.Range("E5").Value = "  +1.86%  "
$ws.Range("D6").Value = "151.59"
$ws.Range("E6").Value = "  -1.10%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E8").Value = "  -3.36%  "
$ws.Range("D9").Value = "6.41"
$ws.Range("E9").Value = "  -3.77%  "
$ws.Range("D10").Value = "0.107"
$ws.Range("E10").Value = "  +2.78%  "
$ws.Range("D11").Value = "0.345"
$ws.Range("E11").Value = "  -0.11%  "
$ws.Range("E12").Value = "  -0.84%  "
$ws.Range("D13").Value = "3.088.42"
$ws.Range("E13").Value = "  +0.92%  "
$ws.Range("D14").Value = "60.512.19"
$ws.Range("E14").Value = "  +0.02%  "
$ws.Range("E15").Value = "  -0.44%  "
$ws.Range("E16").Value = "  -0.29%  "
$ws.Range("D17").Value = "2.622.62"
$ws.Range("E17").Value = "  +0.52%  "
$ws.Range("D18").Value = "4.67"
$ws.Range("E18").Value = "  -1.53%  "
$ws.Range("D19").Value = "347.86"
$ws.Range("E19").Value = "  -2.67%  "
$ws.Range("D20").Value = "10.48"
$ws.Range("E20").Value = "  -1.32%  "
$ws.Range("D21").Value = "6.18"
$ws.Range("E21").Value = "  -0.29%  "
$ws.Range("D22").Value = "0.994"
$ws.Range("E22").Value = "  -0.48%  "
$ws.Range("D23").Value = "61.08"
$ws.Range("E23").Value = "  +0.09%  "
$ws.Range("D25").Value = "0.165"
$ws.Range("E25").Value = "  -0.69%  "
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  +0.29%  "
$ws.Range("E27").Value = "  -0.01%  "
$ws.Range("D28").Value = "7.13"
$ws.Range("E28").Value = "  -1.82%  "
$ws.Range("D29").Value = "0.999"
$ws.Range("E29").Value = "  -0.05%  "
$ws.Range("D30").Value = "6.06"
$ws.Range("E30").Value = "  +2.34%  "
$ws.Range("E31").Value = "  +0.82%  "
$ws.Range("D32").Value = "19.09"
$ws.Range("E32").Value = "  -1.69%  "
$ws.Range("D33").Value = "150.14"
$ws.Range("E33").Value = "  -0.26%  "
$ws.Range("D34").Value = "4.01"
$ws.Range("E34").Value = "  +0.26%  "
$ws.Range("D35").Value = "0.896"
$ws.Range("E35").Value = "  -2.00%  "
$ws.Range("E36").Value = "  -1.98%  "
$ws.Range("D37").Value = "0.882"
$ws.Range("E37").Value = "  +4.73%  "
$ws.Range("D38").Value = "36.69"
$ws.Range("E38").Value = "  +1.36%  "
$ws.Range("E39").Value = "  -1.45%  "
$ws.Range("D40").Value = "3.68"
$ws.Range("E40").Value = "  -1.62%  "
$ws.Range("D41").Value = "291.05"
$ws.Range("E41").Value = "  +1.15%  "
$ws.Range("D42").Value = "0.633"
$ws.Range("E42").Value = "  +2.40%  "
$ws.Range("E43").Value = "  -0.53%  "
$ws.Range("D44").Value = "0.998"
$ws.Range("E44").Value = "  +0.11%  "
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").Value = "19.85"
$ws.Range("E45").Value = "  +1.37%  "
$ws.Range("B46").Value = "Hedera"
$ws.Range("C46").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D46").Value = "0.0553"
$ws.Range("E46").Value = "  -0.45%  "
$ws.Range("E47").Value = "  +0.12%  "
$ws.Range("D48").Value = "4.79"
$ws.Range("E48").Value = "  -3.25%  "
$ws.Range("D49").Value = "10.38"
$ws.Range("E49").Value = "  +0.76%  "
$ws.Range("D50").Value = "18.96"
$ws.Range("E50").Value = "  -1.17%  "
$ws.Range("D51").Value = "1.972.24"
$ws.Range("E51").Value = "  -0.87%  "

# Restore default (Normal) style on the Price cells so no stray text-format
# style lingers on them (matches original workbook formatting).
foreach ($addr in $priceCells) {
    $ws.Range($addr).Style = "Normal"
}
